# Auto-generated edit script: updates cryptos.xlsx D/E (price/volume) columns
# and a few B/C (coin name/link) columns per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $NewValue)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $NewValue
    $r.ClearFormats()
}

Set-TextValue "D2" "64.455.24"
Set-TextValue "E2" "  +0.57%  "
Set-TextValue "D3" "3.141.12"
Set-TextValue "E3" "  -0.14%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "609.16"
Set-TextValue "E5" "  +0.47%  "
Set-TextValue "D6" "143.94"
Set-TextValue "E6" "  -2.10%  "
Set-TextValue "E7" "  -0.04%  "
Set-TextValue "D8" "3.137.89"
Set-TextValue "E8" "  -0.08%  "
Set-TextValue "E9" "  +0.43%  "
Set-TextValue "E10" "  +0.34%  "
Set-TextValue "D11" "5.39"
Set-TextValue "E11" "  -1.10%  "
Set-TextValue "D12" "0.470"
Set-TextValue "E12" "  -0.65%  "
Set-TextValue "D13" "0.0000257"
Set-TextValue "E13" "  +2.87%  "
Set-TextValue "D14" "35.49"
Set-TextValue "E14" "  +0.22%  "
Set-TextValue "D15" "3.658.65"
Set-TextValue "E15" "  -0.08%  "
Set-TextValue "E16" "  +2.58%  "
Set-TextValue "D17" "64.368.72"
Set-TextValue "E17" "  +0.38%  "
Set-TextValue "D18" "3.146.13"
Set-TextValue "E18" "  +0.12%  "
Set-TextValue "E19" "  +0.04%  "
Set-TextValue "D20" "477.43"
Set-TextValue "E20" "  -0.01%  "
Set-TextValue "E21" "  +0.46%  "
Set-TextValue "E22" "  +2.05%  "
Set-TextValue "D23" "7.77"
Set-TextValue "E23" "  -0.07%  "
Set-TextValue "D24" "85.38"
Set-TextValue "E24" "  +2.72%  "
Set-TextValue "D25" "13.45"
Set-TextValue "E25" "  -1.46%  "
Set-TextValue "E26" "  +0.06%  "
Set-TextValue "E27" "  -3.40%  "
Set-TextValue "D28" "8.46"
Set-TextValue "E28" "  +0.33%  "
Set-TextValue "D29" "7.27"
Set-TextValue "E29" "  +7.01%  "
Set-TextValue "E30" "  +2.28%  "
Set-TextValue "E31" "  -5.68%  "
Set-TextValue "B32" "FirstDigitalUSD"
Set-TextValue "C32" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D32" "1.00"
Set-TextValue "E32" "  -0.02%  "
Set-TextValue "B33" "EthereumClassic"
Set-TextValue "C33" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "26.88"
Set-TextValue "E33" "  +3.04%  "
Set-TextValue "D34" "2.64"
Set-TextValue "E34" "  -3.93%  "
Set-TextValue "E35" "  +1.14%  "
Set-TextValue "E36" "  +0.36%  "
Set-TextValue "D37" "0.0₃0760"
Set-TextValue "E37" "  +3.38%  "
Set-TextValue "D38" "52.55"
Set-TextValue "E38" "  -2.01%  "
Set-TextValue "D39" "3.01"
Set-TextValue "E39" "  +2.36%  "
Set-TextValue "D40" "446.70"
Set-TextValue "E40" "  -3.52%  "
Set-TextValue "D41" "0.0395"
Set-TextValue "E41" "  +0.35%  "
Set-TextValue "D42" "0.120"
Set-TextValue "E42" "  +1.23%  "
Set-TextValue "D43" "8.28"
Set-TextValue "E43" "  -1.51%  "
Set-TextValue "D44" "2.898.10"
Set-TextValue "E44" "  +1.83%  "
Set-TextValue "D45" "0.262"
Set-TextValue "E45" "  -1.11%  "
Set-TextValue "D46" "2.23"
Set-TextValue "E46" "  -1.42%  "
Set-TextValue "D47" "2.41"
Set-TextValue "E47" "  +3.48%  "
Set-TextValue "D48" "26.33"
Set-TextValue "E48" "  -0.42%  "
Set-TextValue "E50" "  -0.63%  "
Set-TextValue "B51" "Monero"
Set-TextValue "C51" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D51" "120.11"
Set-TextValue "E51" "  +0.73%  "
